# "commit Today assignment and teamReport"
# Fill in today's attendance row (row 7, 2023-08-09) and note down the two
# people who were absent, same as the rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 attendance: everyone Present except Pravin Mane (E7) and
# Prathmesh Patil (F7), who are Absent.
$ws.Range("B7").Value = "Present"
$ws.Range("C7").Value = "Present"
$ws.Range("D7").Value = "Present"
$ws.Range("E7").Value = "Absent"
$ws.Range("F7").Value = "Absent"
$ws.Range("G7").Value = "Present"
$ws.Range("H7").Value = "Present"
$ws.Range("I7").Value = "Absent"

# Explain the two absences with comments, same style as the existing notes
# (author "HP:" lead-in followed by the reason).
$commentE7 = $ws.Range("E7").AddComment("HP:" + [char]10 + "Due to personal reason unable to join the session.")
$commentF7 = $ws.Range("F7").AddComment("HP:" + [char]10 + "he has not join the meeting because he has left from office very late.")

# Leave the selection where it landed after filling the row in.
$null = $ws.Range("E10").Select()
